$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DigitalCollection")
$ws2 = $wb.Worksheets.Item("SpaceDetails")

# ----------------------------------------------------------------------
# DigitalCollection sheet
# ----------------------------------------------------------------------

# Fix typo: "Art in Christian Tradition" -> "Art in the Christian Tradition"
$ws1.Range("A7").Value = "Art in the Christian Tradition"

# New row 14: History of Art Image Repository
$ws1.Range("A14").Value = "History of Art Image Repository"
$ws1.Range("B14").Value = "https://library-artstor-org.proxy.library.vanderbilt.edu/#/collection/37831/"
$ws1.Range("C14").Value = 104423
$ws1.Range("C14").NumberFormat = "#,##0"
$ws1.Range("D14").Value = "from Cliff"

# Widen column A to fit the longer collection names
$ws1.Columns.Item(1).ColumnWidth = 31.666666666666668

# ----------------------------------------------------------------------
# SpaceDetails sheet
# ----------------------------------------------------------------------

# Fix typo: "Annett" -> "Annette"
$ws2.Range("A4").Value = "Annette & Irwin Eskind Family Biomedical Library and Learning Center"

# Update counts
$ws2.Range("C3").Value = 198
$ws2.Range("C5").Value = 418
$ws2.Range("B6").Value = 12600
$ws2.Range("E6").Value = 7
$ws2.Range("C7").Value = 86
$ws2.Range("C9").Value = 326
$ws2.Range("C10").Value = 34
$ws2.Range("D10").Value = 1

# New column A width so the longer library names fit
$ws2.Columns.Item(1).ColumnWidth = 37.333333333333336

# ----------------------------------------------------------------------
# Sheet views / selection: DigitalCollection becomes the active tab,
# selection on DigitalCollection!D14 and SpaceDetails!A13
# ----------------------------------------------------------------------

[void]$ws2.Activate()
[void]$ws2.Range("A13").Select()
[void]$ws1.Activate()
[void]$ws1.Range("D14").Select()
